$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich-text cells) ---
# "Volume 30   Number  34" -> "...35"
$ws.Range("A8").Value = "Volume 30   Number  35"
# "Report Covering the Week  8/21/2023  Through  8/27/2023" -> new week dates
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Row 15 (Rape): 28-day % chg column only ---
$ws.Range("L15").Value = 100

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = 181.25
$ws.Range("M16").Value = -19.642857142857
$ws.Range("N16").Value = -79.45205479452

# --- Row 17 (Fel. Assault) ---
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 43.478260869565
$ws.Range("M17").Value = 78.378378378378
$ws.Range("N17").Value = -15.384615384615

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 15.78947368421
$ws.Range("I18").Value = 210
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = 20.689655172413
$ws.Range("L18").Value = 43.835616438356
$ws.Range("M18").Value = 26.506024096385
$ws.Range("N18").Value = -70.128022759601

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -38.095238095238
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -18.965517241379
$ws.Range("I19").Value = 445
$ws.Range("J19").Value = 421
$ws.Range("K19").Value = 5.700712589073
$ws.Range("L19").Value = 79.435483870967
$ws.Range("M19").Value = 80.894308943089
$ws.Range("N19").Value = 26.420454545454

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 54.545454545454
$ws.Range("I20").Value = 112
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = 138.297872340426
$ws.Range("M20").Value = 19.148936170212
$ws.Range("N20").Value = -94.759007955077

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -16.129032258064
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = 0.980392156862
$ws.Range("I21").Value = 888
$ws.Range("J21").Value = 787
$ws.Range("K21").Value = 12.833545108005
$ws.Range("L21").Value = 75.147928994082
$ws.Range("M21").Value = 47.263681592039
$ws.Range("N21").Value = -74.62132037725

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 36.363636363636
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = -5.66037735849
$ws.Range("I24").Value = 392
$ws.Range("J24").Value = 517
$ws.Range("K24").Value = -24.177949709864
$ws.Range("L24").Value = 7.397260273972
$ws.Range("M24").Value = 26.04501607717

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 3
$ws.Range("E25").Value = -25
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -5.882352941176
$ws.Range("I25").Value = 143
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = -10.625
$ws.Range("L25").Value = 55.434782608695
$ws.Range("M25").Value = 24.347826086956

# --- Row 26 (UCR Rape*): 28-day % chg column only ---
$ws.Range("L26").Value = 0

# --- Row 27 (Other Sex Crimes): 28-day % chg column only ---
$ws.Range("L27").Value = -25
